# Merge split "N. " + "title" runs into a single run for sections 2-8.
$d = $word.ActiveDocument

$titles = @(
    '2. Truy cập và thao tác phần tử',
    '3. Sắp xếp mảng',
    '4. Tìm kiếm và lọc',
    '5. Tìm kiếm và lọc',
    '6. Lấy key và value',
    '7. Ánh xạ và xử lý',
    '8. Các hàm tiện ích khác'
)

foreach ($title in $titles) {
    $null = $d.Content.Find.Execute($title, $true, $false, $false, $false, $false, $true, 1, $false, $title, 2)
}

# Insert the new block of paragraphs (array_key_exists ... array_intersect_assos)
# right before the trailing empty paragraph, then remove that now-redundant
# trailing empty paragraph.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint = $lastPara.Range
$insertPoint.Collapse(1)

$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve">array_key_exists($key, $array): </w:t>
  </w:r>
  <w:r>
    <w:t>kiểm tra khóa $value có tồn tại trong mảng không, nếu có trả về true</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve">implode($str, $array): </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">chuyển các giá trị của mảng thành một chuỗi bao gồm các phần tử cách nhau bởi kí tự </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>$str</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve">explore($delimiter, $str): </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">chuyển chuỗi thành mảng, tách chuỗi dựa vào </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>$delimiter</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:t>9. Truy xuất phần tử của mảng</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>c</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve">urrent($array): </w:t>
  </w:r>
  <w:r>
    <w:t>truy xuất phần tử hiện tại của mảng</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>e</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve">nd($array): </w:t>
  </w:r>
  <w:r>
    <w:t>truy xuất phần tử cuối của mảng</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>n</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve">ext($array): </w:t>
  </w:r>
  <w:r>
    <w:t>truy xuất phần tử sau phần tử hiện tại</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>p</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve">rev($array): </w:t>
  </w:r>
  <w:r>
    <w:t>truy xuất phần tử trước phần tử hiện tại</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>r</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve">eset($array): </w:t>
  </w:r>
  <w:r>
    <w:t>quay về vị trí phần tử đầu tiên trong mảng.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:t>10. Các trường hợp so sánh giữa hai mảng</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:t>TH1: So sánh khác nhau (array_diff):</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>a</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>rray.diff($a1, $a2,…):</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>trả về những phần tử của mảng thứ nhất mà không tồn tại trong các mảng còn lại.</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> (aray_diff_key – So sánh dựa vào key)</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t>TH2: So sánh giống nhau</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve">array_intersect($a1, $a2,…): </w:t>
  </w:r>
  <w:r>
    <w:t>trả về một mảng bao gồm các phần tử giống nhau về giá trị giữa 2 mảng $a1 và $a2</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve">array_intersect_key($a1, $a2,…): </w:t>
  </w:r>
  <w:r>
    <w:t>trả về một mảng các phần tử giống nhau về 2 khóa giữa 2 mảng</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve">array_intersect_assos($a1, $a2,…): </w:t>
  </w:r>
  <w:r>
    <w:t>trả về một mảng bao gồm các phần tử giống nhau về khóa và giá trị giữa 2 mảng</w:t>
  </w:r>
</w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertPoint.InsertXML($xml)

$markerRange = $d.Content
$null = $markerRange.Find.Execute('trả về một mảng bao gồm các phần tử giống nhau về khóa và giá trị giữa 2 mảng')
$tailRange = $d.Range($markerRange.End, $d.StoryRanges.Item(1).End)
$tailRange.Delete()
